$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (duplicate of row 2) and row 4 held identical machine data; row 4 also
# carried a handful of blank-but-present cells (C/E/F/Z) from a wider paste.
# Drop row 3 so row 4 (with those blank cells) shifts up into its place,
# collapsing the two duplicate rows into one while keeping the extra blanks.
$ws.Rows.Item(3).Delete()

# "N/A" -> "Nenhum" for the "Troca ou Upgrade" column on both remaining data rows.
$ws.Range("T2").Value = "Nenhum"

# New priority note for the (now single) data row.
$ws.Range("U3").Value = "Não será trocada"
